$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 99:100 (existing rows 99+ shift down to 101+)
$ws.Rows("99:100").Insert()

# ---- Fill in the new row 99 ----
$ws.Range("A99").Value = 7
$ws.Range("B99").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C99").Value = "Ñuble"
$ws.Range("D99").Value = 44466
$ws.Range("E99").Value = 16
$ws.Range("F99").Value = 100112002
$ws.Range("G99").Value = "Pimiento"
$ws.Range("H99").Value = "Zafiro rojo"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 60
$ws.Range("K99").Value = 41000
$ws.Range("L99").Value = 42000
$ws.Range("M99").Value = 41500
$ws.Range("N99").Value = "$/caja 15 kilos"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 2767
$ws.Range("Q99").Value = 15
$ws.Range("R99").Value = "Hortaliza"

# ---- Fill in the new row 100 ----
$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 44466
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112002
$ws.Range("G100").Value = "Pimiento"
$ws.Range("H100").Value = "Zafiro verde"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 100
$ws.Range("K100").Value = 34000
$ws.Range("L100").Value = 35000
$ws.Range("M100").Value = 34500
$ws.Range("N100").Value = "$/caja 15 kilos"
$ws.Range("O100").Value = "Región de Arica y Parinacota"
$ws.Range("P100").Value = 2300
$ws.Range("Q100").Value = 15
$ws.Range("R100").Value = "Hortaliza"
